# Refined metadata to be additional tab
#
# 1. Update the "panel_query_time"-equivalent timestamps (column F) on the
#    existing "data" sheet.
# 2. Add a new "metadata" worksheet right after "data", with a header row
#    (bold / bordered, matching the "data" sheet header style) and one row
#    of panel metadata.

$wb = $excel.ActiveWorkbook
$ds = $wb.Worksheets.Item("data")

# --- 1. Refresh timestamps on the "data" sheet -----------------------------
$ds.Range("F2").Value = "2021-10-05 14:22:20.163145"
$ds.Range("F3").Value = "2021-10-05 14:22:20.163153"
$ds.Range("F4").Value = "2021-10-05 14:22:20.163155"
$ds.Range("F5").Value = "2021-10-05 14:22:20.163157"

# --- 2. Create the new "metadata" sheet, placed after "data" ---------------
$meta = $wb.Worksheets.Add($null, $ds)
$meta.Name = "metadata"

# Header row (B1:G1)
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Copy the header formatting (bold, border, centered) from the "data" sheet
# header so the new header reuses the exact same style.
$ds.Range("B1:F1").Copy()
$meta.Range("B1:F1").PasteSpecial(-4122)
$ds.Range("B1").Copy()
$meta.Range("G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data row (A2:G2)
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Primary pigmented nodular adrenocortical disease"
$meta.Range("C2").Value = 566
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.7"
$meta.Range("E2").Value = "2021-03-02T16:14:35.738760Z"
$meta.Range("F2").Value = "2021-10-05 14:22:20.159135"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/566/?format=json"

# Copy the A-column style (bold/centered/bordered, used for the numeric
# index column) from the "data" sheet onto the new sheet's A2 cell.
$ds.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$meta.Range("A2").Value = 0
